$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'35.516.39"
$ws.Range("E2").Value = "  +1.43%  "

# Row 3
$ws.Range("D3").Value = "'1.877.17"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.37%  "

# Row 5
$ws.Range("D5").Value = "'242.73"
$ws.Range("E5").Value = "  +4.19%  "

# Row 6
$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "  +2.32%  "

# Row 7
$ws.Range("E7").Value = "  +0.52%  "

# Row 8
$ws.Range("D8").Value = "'43.45"
$ws.Range("E8").Value = "  +6.46%  "

# Row 9
$ws.Range("E9").Value = "  +0.50%  "

# Row 10
$ws.Range("E10").Value = "  +1.44%  "

# Row 11
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("D12").Value = "'2.145.22"
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("D13").Value = "'11.97"
$ws.Range("E13").Value = "  +4.88%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.888.11"
$ws.Range("E14").Value = "  +2.06%  "

# Row 15
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.686"
$ws.Range("E15").Value = "  +1.54%  "

# Row 16
$ws.Range("E16").Value = "  +2.21%  "

# Row 17
$ws.Range("D17").Value = "'35.532.83"
$ws.Range("E17").Value = "  +1.35%  "

# Row 18
$ws.Range("D18").Value = "'71.12"
$ws.Range("E18").Value = "  +1.44%  "

# Row 19
$ws.Range("D19").Value = "'0.0₃0806"
$ws.Range("E19").Value = "  +1.66%  "

# Row 20
$ws.Range("D20").Value = "'243.13"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21
$ws.Range("D21").Value = "'12.36"
$ws.Range("E21").Value = "  +0.67%  "

# Row 22
$ws.Range("D22").Value = "'4.83"
$ws.Range("E22").Value = "  +1.54%  "

# Row 23
$ws.Range("E23").Value = "  +0.53%  "

# Row 24
$ws.Range("D24").Value = "'2.29"
$ws.Range("E24").Value = "  +1.45%  "

# Row 25
$ws.Range("D25").Value = "'171.45"
$ws.Range("E25").Value = "  -0.71%  "

# Row 26
$ws.Range("E26").Value = "  +27.66%  "

# Row 27
$ws.Range("D27").Value = "'8.27"
$ws.Range("E27").Value = "  +5.32%  "

# Row 28
$ws.Range("D28").Value = "'17.85"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("E30").Value = "  +1.70%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.06"
$ws.Range("E31").Value = "  +2.71%  "

# Row 32
$ws.Range("B32").Value = "BinanceUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D32").Value = "'1.02"
$ws.Range("E32").Value = "  +0.59%  "

# Row 33
$ws.Range("D33").Value = "'0.925"
$ws.Range("E33").Value = "  +22.33%  "

# Row 34
$ws.Range("D34").Value = "'4.09"
$ws.Range("E34").Value = "  +2.74%  "

# Row 35
$ws.Range("D35").Value = "'1.78"
$ws.Range("E35").Value = "  +11.70%  "

# Row 36
$ws.Range("E36").Value = "  +5.16%  "

# Row 37
$ws.Range("D37").Value = "'1.35"
$ws.Range("E37").Value = "  +11.15%  "

# Row 38
$ws.Range("E38").Value = "  +2.02%  "

# Row 39
$ws.Range("E39").Value = "  +4.12%  "

# Row 40
$ws.Range("D40").Value = "'90.00"
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("D41").Value = "'1.355.92"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42
$ws.Range("D42").Value = "'15.22"
$ws.Range("E42").Value = "  +3.85%  "

# Row 43
$ws.Range("E43").Value = "  +11.09%  "

# Row 44
$ws.Range("E44").Value = "  +3.27%  "

# Row 45
$ws.Range("E45").Value = "  +0.46%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'6.72"
$ws.Range("E46").Value = "  +6.15%  "

# Row 47
$ws.Range("B47").Value = "Gas"
$ws.Range("C47").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D47").Value = "'12.44"
$ws.Range("E47").Value = "  +40.46%  "

# Row 48
$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").Value = "'2.74"
$ws.Range("E48").Value = "  -0.92%  "

# Row 49
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "'44.81"
$ws.Range("E49").Value = "  +32.03%  "

# Row 50
$ws.Range("D50").Value = "'2.064.64"
$ws.Range("E50").Value = "  +1.31%  "

# Row 51
$ws.Range("E51").Value = "  +2.41%  "
